$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.775.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +8.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.948.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.22"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4784"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4136"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.96"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08267"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.047"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.932.47"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.173"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.446"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06712"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.728.12"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.624"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.89%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.169.01"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.70"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.198"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.651"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.027"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09655"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.474"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +11.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.684"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.494"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.11%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06250"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.37%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02314"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.716"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.12%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6111"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.75"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1910"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.276"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5720"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.45%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.54"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.322"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +27.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07463"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +13.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.12%  "
